# Rebuild the "rand_var" table: transpose the old 10-row (class interval / X)
# layout into a 3-row (lower bound / upper bound / X) layout spanning many
# columns, with updated labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so no stale cells/styles from the old 10-row
# layout linger outside the new A2:K4 range.
$ws.Cells.Clear()

# --- Row 2: lower class bounds -------------------------------------------------
$ws.Cells.Item(2, 1).Value = "Weights (Lower class bound)"
$lowerBounds = 31, 36, 41, 46, 51, 56, 61, 66, 71, 76
for ($i = 0; $i -lt $lowerBounds.Length; $i++) {
    $ws.Cells.Item(2, 2 + $i).Value = $lowerBounds[$i]
}

# --- Row 3: upper class bounds -------------------------------------------------
$ws.Cells.Item(3, 1).Value = "Weights (Upper class bound)"
$upperBounds = 35, 40, 45, 50, 55, 60, 65, 70, 75, 80
for ($i = 0; $i -lt $upperBounds.Length; $i++) {
    $ws.Cells.Item(3, 2 + $i).Value = $upperBounds[$i]
}

# --- Row 4: X (random variable) -------------------------------------------------
$ws.Cells.Item(4, 1).Value = "X (Random Variable)"
$xValues = 1, 2, 3, 4, 5, 6, 7, 8, 9
for ($i = 0; $i -lt $xValues.Length; $i++) {
    $ws.Cells.Item(4, 2 + $i).Value = $xValues[$i]
}

# --- Formatting -----------------------------------------------------------------
# Row labels (column A) are bold.
$ws.Range("A2:A4").Font.Bold = $true

# The old "center" alignment style becomes "right" alignment, still present on
# column C (and the lingering B2 cell that inherited the old row-2 style).
$ws.Range("B2").HorizontalAlignment = -4152
$ws.Range("C2:C4").HorizontalAlignment = -4152

# Column widths: label column is widened, column C keeps its (now right
# aligned) default width.
$ws.Columns("A").ColumnWidth = 24.14

# --- Selection / view ------------------------------------------------------------
[void]$ws.Range("K3").Select()
